$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in the "S"/"T" columns (extra points / remark) for the listed students.
$ws.Range("S5").Value = 5
$ws.Range("T5").Value = "автоматом"

$ws.Range("S6").Value = 3
$ws.Range("T6").Value = "авансом"

$ws.Range("S7").Value = 3
$ws.Range("T7").Value = "авансом"

$ws.Range("S8").Value = "отчислена"

$ws.Range("S18").Value = 3
$ws.Range("T18").Value = "авансом"

$ws.Range("S22").Value = 4
$ws.Range("T22").Value = "авансом"

$ws.Range("S23").Value = 3
$ws.Range("T23").Value = "авансом"

$ws.Range("S31").Value = 3
$ws.Range("T31").Value = "авансом"

$ws.Range("S33").Value = 4
$ws.Range("T33").Value = "авансом"

# Update the selection to match the latest edit position.
$ws.Range("S9").Select()
